$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder fields: "4/3/15" -> "9/4/15"
#    (slide master + every custom layout's "Date Placeholder *" shape)
# ---------------------------------------------------------------------------
function Set-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "9/4/15"
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Set-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# Slide 1 edits
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# 2. "Obtained from NOAA/NGDC" -> "Obtained from NOAA/NCEI"
$shpNgdc = $s.Shapes.Item(6)
$found = $shpNgdc.TextFrame.TextRange.Find("NGDC")
$found.Text = "NCEI"

# 3. "Facular Brightening Function (PX)" -> "Facular Brightening Function F(t)"
#    with F and t in italics
$shpFacular = $s.Shapes.Item(11)
$trFacular = $shpFacular.TextFrame.TextRange
$trFacular.Text = "Facular Brightening Function F(t)"
$trFacular.Characters(30, 1).Font.Italic = $true   # "F"
$trFacular.Characters(32, 1).Font.Italic = $true   # "t"

# 4. "Sunspot Blocking Function (PS)" -> "Sunspot Blocking Function S(t)"
#    with S and t in italics
$shpSunspot = $s.Shapes.Item(12)
$trSunspot = $shpSunspot.TextFrame.TextRange
$trSunspot.Text = "Sunspot Blocking Function S(t)"
$trSunspot.Characters(27, 1).Font.Italic = $true   # "S"
$trSunspot.Characters(29, 1).Font.Italic = $true   # "t"

# 5. "Compute Model Inputs" + ":" runs -> single run "Compute Model Inputs:"
$shpCompute = $s.Shapes.Item(16)
$trCompute = $shpCompute.TextFrame.TextRange
$trCompute.Text = ""
$trCompute.Text = "Compute Model Inputs:"
$trCompute.Font.Underline = $true

# 6. "Generate Output Files" + ":" runs -> single run "Generate Output Files:"
$shpGenerate = $s.Shapes.Item(18)
$trGenerate = $shpGenerate.TextFrame.TextRange
$trGenerate.Text = ""
$trGenerate.Text = "Generate Output Files:"
$trGenerate.Font.Underline = $true

# 7. "Scaling Factors convert PX and PS to Irradiance Changes"
#    -> "Scaling Factors convert F(t) and S(t) to Irradiance Changes"
#    with F, t, S, t in italics. Only the first paragraph of this text box
#    changes; the remaining paragraphs are left untouched.
$shpScaling = $s.Shapes.Item(25)
$trScaling = $shpScaling.TextFrame.TextRange
$oldHeading = $trScaling.Find("Scaling Factors convert PX and PS to Irradiance Changes")
$oldHeading.Text = "Scaling Factors convert F(t) and S(t) to Irradiance Changes"
$trScaling.Characters(25, 1).Font.Italic = $true   # "F"
$trScaling.Characters(27, 1).Font.Italic = $true   # "t"
$trScaling.Characters(34, 1).Font.Italic = $true   # "S"
$trScaling.Characters(36, 1).Font.Italic = $true   # "t"
